$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '70.869.83'
$ws.Range("E2").Value = '  +6.78%  '

$ws.Range("D3").Value = '3.635.01'
$ws.Range("E3").Value = '  +6.30%  '

$ws.Range("E4").Value = '  +0.02%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '592.94'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +4.68%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '192.24'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +7.79%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.653'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +3.26%  '

$ws.Range("D8").Value = '3.625.63'
$ws.Range("E8").Value = '  +6.31%  '

$ws.Range("E9").Value = '  +0.08%  '

$ws.Range("E10").Value = '  +2.39%  '

$ws.Range("E11").Value = '  +4.37%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '58.17'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +6.71%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000298'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +6.15%  '

$ws.Range("E14").Value = '  +5.26%  '

$ws.Range("D15").Value = '4.214.65'
$ws.Range("E15").Value = '  +7.13%  '

$ws.Range("D16").Value = '3.632.13'
$ws.Range("E16").Value = '  +6.67%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '19.41'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +5.90%  '

$ws.Range("D18").Value = '70.779.77'
$ws.Range("E18").Value = '  +6.98%  '

$ws.Range("E19").Value = '  +5.40%  '

$ws.Range("E20").Value = '  +0.68%  '

$ws.Range("E21").Value = '  +4.97%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '496.72'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +6.33%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.42'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +9.46%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '17.26'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +16.76%  '

$ws.Range("E25").Value = '  +8.89%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '91.19'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +1.24%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '3.13'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +6.25%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '11.30'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +4.86%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '9.48'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +7.34%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '32.44'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +3.31%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '7.58'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +12.99%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '12.24'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +6.03%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '619.08'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +6.22%  '

$ws.Range("E34").Value = '  +8.20%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '65.32'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +4.47%  '

$ws.Range("B36").Value = 'TheGraph'
$ws.Range("C36").Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.417'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +9.14%  '

$ws.Range("B37").Value = 'PEPE'
$ws.Range("C37").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D37").Value = '0.0₃0832'
$ws.Range("E37").Value = '  +8.97%  '

$ws.Range("E38").Value = '  +2.98%  '

$ws.Range("E39").Value = '  +4.61%  '

$ws.Range("E40").Value = '  -0.03%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '3.67'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +2.59%  '

$ws.Range("D42").Value = '3.333.80'
$ws.Range("E42").Value = '  +6.48%  '

$ws.Range("E43").Value = '  +6.71%  '

$ws.Range("E44").Value = '  +6.54%  '

$ws.Range("E45").Value = '  +8.49%  '

$ws.Range("E46").Value = '  +5.46%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.139'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +3.04%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '9.22'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +7.41%  '

$ws.Range("E49").Value = '  +3.81%  '

$ws.Range("E50").Value = '  +4.70%  '

$ws.Range("B51").Value = 'FirstDigitalUSD'
$ws.Range("C51").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.999'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.07%  '

